$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows at 207-210 (shifts old rows 207-255 down to 211-259)
$ws.Range("A207:A210").EntireRow.Insert()

# Row 207
$ws.Cells.Item(207, 1).Value = 10
$ws.Cells.Item(207, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(207, 3).Value = 'La Araucanía'
$ws.Cells.Item(207, 4).Value = 44900
$ws.Cells.Item(207, 5).Value = 9
$ws.Cells.Item(207, 6).Value = 'Fruta'
$ws.Cells.Item(207, 7).Value = 100103
$ws.Cells.Item(207, 8).Value = 'Frutos de hueso (carozo)'
$ws.Cells.Item(207, 9).Value = 100103001
$ws.Cells.Item(207, 10).Value = 'Cereza'
$ws.Cells.Item(207, 11).Value = 'Brooks'
$ws.Cells.Item(207, 12).Value = 'Primera'
$ws.Cells.Item(207, 13).Value = 1000
$ws.Cells.Item(207, 14).Value = 5500
$ws.Cells.Item(207, 15).Value = 6000
$ws.Cells.Item(207, 16).Value = 5750
$ws.Cells.Item(207, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(207, 18).Value = 'Región del Maule'
$ws.Cells.Item(207, 19).Value = 575
$ws.Cells.Item(207, 20).Value = 10

# Row 208
$ws.Cells.Item(208, 1).Value = 10
$ws.Cells.Item(208, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(208, 3).Value = 'La Araucanía'
$ws.Cells.Item(208, 4).Value = 44900
$ws.Cells.Item(208, 5).Value = 9
$ws.Cells.Item(208, 6).Value = 'Fruta'
$ws.Cells.Item(208, 7).Value = 100103
$ws.Cells.Item(208, 8).Value = 'Frutos de hueso (carozo)'
$ws.Cells.Item(208, 9).Value = 100103001
$ws.Cells.Item(208, 10).Value = 'Cereza'
$ws.Cells.Item(208, 11).Value = 'Brooks'
$ws.Cells.Item(208, 12).Value = 'Primera'
$ws.Cells.Item(208, 13).Value = 1800
$ws.Cells.Item(208, 14).Value = 600
$ws.Cells.Item(208, 15).Value = 700
$ws.Cells.Item(208, 16).Value = 644
$ws.Cells.Item(208, 17).Value = '$/kilo (en caja de 15 kilos)'
$ws.Cells.Item(208, 18).Value = 'Región del Maule'
$ws.Cells.Item(208, 19).Value = 644
$ws.Cells.Item(208, 20).Value = 1

# Row 209
$ws.Cells.Item(209, 1).Value = 10
$ws.Cells.Item(209, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(209, 3).Value = 'La Araucanía'
$ws.Cells.Item(209, 4).Value = 44900
$ws.Cells.Item(209, 5).Value = 9
$ws.Cells.Item(209, 6).Value = 'Fruta'
$ws.Cells.Item(209, 7).Value = 100103
$ws.Cells.Item(209, 8).Value = 'Frutos de hueso (carozo)'
$ws.Cells.Item(209, 9).Value = 100103001
$ws.Cells.Item(209, 10).Value = 'Cereza'
$ws.Cells.Item(209, 11).Value = 'Rainier'
$ws.Cells.Item(209, 12).Value = 'Especial'
$ws.Cells.Item(209, 13).Value = 500
$ws.Cells.Item(209, 14).Value = 12000
$ws.Cells.Item(209, 15).Value = 12000
$ws.Cells.Item(209, 16).Value = 12000
$ws.Cells.Item(209, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(209, 18).Value = 'Región del Maule'
$ws.Cells.Item(209, 19).Value = 1200
$ws.Cells.Item(209, 20).Value = 10

# Row 210
$ws.Cells.Item(210, 1).Value = 10
$ws.Cells.Item(210, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(210, 3).Value = 'La Araucanía'
$ws.Cells.Item(210, 4).Value = 44900
$ws.Cells.Item(210, 5).Value = 9
$ws.Cells.Item(210, 6).Value = 'Fruta'
$ws.Cells.Item(210, 7).Value = 100103
$ws.Cells.Item(210, 8).Value = 'Frutos de hueso (carozo)'
$ws.Cells.Item(210, 9).Value = 100103001
$ws.Cells.Item(210, 10).Value = 'Cereza'
$ws.Cells.Item(210, 11).Value = 'Rainier'
$ws.Cells.Item(210, 12).Value = 'Especial'
$ws.Cells.Item(210, 13).Value = 700
$ws.Cells.Item(210, 14).Value = 1300
$ws.Cells.Item(210, 15).Value = 1500
$ws.Cells.Item(210, 16).Value = 1386
$ws.Cells.Item(210, 17).Value = '$/kilo (en caja de 15 kilos)'
$ws.Cells.Item(210, 18).Value = 'Región del Maule'
$ws.Cells.Item(210, 19).Value = 1386
$ws.Cells.Item(210, 20).Value = 1
